$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.045.76'
$ws.Range("D3").Value = '3.218.10'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '396.44'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.77'
$ws.Range("E6").Value = '  +7.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.552'
$ws.Range("E7").Value = '  +2.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("E9").Value = '  +5.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.11'
$ws.Range("E10").Value = '  +6.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0916'
$ws.Range("E11").Value = '  +7.07%  '
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("D13").Value = '3.729.00'
$ws.Range("E13").Value = '  +4.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.05'
$ws.Range("E14").Value = '  +3.89%  '
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("D16").Value = '3.194.04'
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("E17").Value = '  +5.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.78'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").Value = '55.910.20'
$ws.Range("E19").Value = '  +8.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("E21").Value = '  +7.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.92'
$ws.Range("E22").Value = '  +3.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '297.70'
$ws.Range("E23").Value = '  +12.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.33'
$ws.Range("E24").Value = '  +7.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.20'
$ws.Range("E25").Value = '  +1.85%  '
$ws.Range("E26").Value = '  +1.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.05'
$ws.Range("E27").Value = '  +2.78%  '
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("E29").Value = '  +4.24%  '
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("E31").Value = '  +4.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.10'
$ws.Range("E32").Value = '  +6.75%  '
$ws.Range("E33").Value = '  +3.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.17'
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.12'
$ws.Range("E37").Value = '  +25.42%  '
$ws.Range("E38").Value = '  +5.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '136.52'
$ws.Range("E40").Value = '  +5.48%  '
$ws.Range("E41").Value = '  +5.23%  '
$ws.Range("E42").Value = '  +3.09%  '
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("E44").Value = '  +3.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.284'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.16'
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("E47").Value = '  +52.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.11'
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").Value = '2.118.35'
$ws.Range("E50").Value = '  +2.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0360'
$ws.Range("E51").Value = '  +10.30%  '
